# WS02-schedule-wishes: refresh "Name" column with alphabetically sorted,
# corrected names (last-name initial / typo fixes) and tidy up the manually
# added rows at the bottom of the sheet, per commit message:
# "All names now have the first letter of their last name added."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Name") values, row by row (row 1 is the header).
$ws.Range("E2").Value  = "Anne-Mai Pall"
$ws.Range("E4").Value  = "Buse Erdem"
$ws.Range("E5").Value  = "Daniell Lepp"
$ws.Range("E6").Value  = "Deniz Gülmez"
$ws.Range("E7").Value  = "Fred Kasemaa"
$ws.Range("E8").Value  = "Helena Jallai"
$ws.Range("E9").Value  = "Helena Mandel"
$ws.Range("E10").Value = "Joel Kikerpill"
$ws.Range("E11").Value = "Joonas Kari Kuusik"
$ws.Range("E12").Value = "Kevin Vahtra"
$ws.Range("E13").Value = "SMJuulia Kaas"
$ws.Range("E14").Value = "Kristjan Kaljurand"
$ws.Range("E15").Value = "Madis Klesment"
$ws.Range("E16").Value = "Martin Merisalu"
$ws.Range("E17").Value = "Rasmus Rahnu"
$ws.Range("E18").Value = "Robert Matjus"
$ws.Range("E19").Value = "Shawn Michael Rains"
$ws.Range("E20").Value = "Siim Lillemets"
$ws.Range("E21").Value = "Siim-Kaarel Kabel"
$ws.Range("E23").Value = "Siim-Sander Sägi"
$ws.Range("E24").Value = "Sven-Ervin Paap"

# Comment cell describing the manually added rows below the survey answers.
$ws.Range("C26").Value = "Siia panin nimesid manuaalselt juurde"

# Remaining manually-added names, continuing the alphabetical ordering.
$ws.Range("E26").Value = "Tanel Madisson"
$ws.Range("E27").Value = "Timo Kirpu"
$ws.Range("E28").Value = "Triinu Saks"
$ws.Range("E29").Value = "Uku Jaan Leppik"
$ws.Range("E30").Value = "Urmo Olesk"

# Last manually-added row is no longer used.
$ws.Range("E31").ClearContents()

# Table got rebuilt/renamed during the rework of the parser.
$lo = $ws.ListObjects.Item(1)
$lo.Name = "Table13"

# Leave selection where the edits ended, matching the manual editing flow.
$ws.Range("B32").Select()
